# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 09:19"

# Simple numeric refreshes (no country re-ordering involved)
# Estados Unidos (row 4)
$ws.Range("B4").Value = 2263749
$ws.Range("C4").Value = 98
$ws.Range("D4").Value = 931076
$ws.Range("E4").Value = 1211985

# Ucrania (row 38)
$ws.Range("B38").Value = 34984
$ws.Range("C38").Value = 921
$ws.Range("D38").Value = 16033
$ws.Range("E38").Value = 17966
$ws.Range("G38").Value = 19
$ws.Range("H38").Value = 985

# Armenia (row 52)
$ws.Range("B52").Value = 19157
$ws.Range("C52").Value = 459
$ws.Range("D52").Value = 8266
$ws.Range("E52").Value = 10572
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 319

# Chequia (row 67)
$ws.Range("B67").Value = 10283
$ws.Range("C67").Value = 3
$ws.Range("D67").Value = 7446
$ws.Range("E67").Value = 2503

# Kenia / El Salvador swap places (rows 85-86), with El Salvador's numbers refreshed
$ws.Range("A85").Value = "El Salvador"
$ws.Range("B85").Value = 4329
$ws.Range("C85").Value = 129
$ws.Range("D85").Value = 2310
$ws.Range("E85").Value = 1933
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 86

$ws.Range("A86").Value = "Kenia"
$ws.Range("B86").Value = 4257
$ws.Range("C86").Value = 0
$ws.Range("D86").Value = 1459
$ws.Range("E86").Value = 2681
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 117

# Hungria (row 88)
$ws.Range("B88").Value = 4081
$ws.Range("C88").Value = 2
$ws.Range("D88").Value = 2581
$ws.Range("E88").Value = 932

# Albania / Lituania swap places (rows 110-111), with Lituania's numbers refreshed
$ws.Range("A110").Value = "Lituania"
$ws.Range("B110").Value = 1792
$ws.Range("C110").Value = 8
$ws.Range("D110").Value = 1462
$ws.Range("E110").Value = 254
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 76

$ws.Range("A111").Value = "Albania"
$ws.Range("B111").Value = 1788
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 1086
$ws.Range("E111").Value = 663
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 39

# Letonia (row 125)
$ws.Range("B125").Value = 1110
$ws.Range("C125").Value = 2
$ws.Range("E125").Value = 177

# Islas Turcas y Caicos / Santa Sede swap places (rows 208-209), pure swap
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 12
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("B209").Value = 12
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

# Islas Virgenes Britanicas / Papua Nueva Guinea swap places (rows 213-214), pure swap
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 8
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
